# Update TPM-derived NATMI ligand-receptor metrics with newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 3.281109666666667
$ws.Range("N2").Value = 9.843329000000001
$ws.Range("O2").Value = 0.2779739143628921
$ws.Range("P2").Value = 0.2779739143628921
$ws.Range("Q2").Value = 0.07725591450811113
$ws.Range("R2").Value = 0.6953032305730001
$ws.Range("S2").Value = 0.0007826509535548598
$ws.Range("T2").Value = 0.0007826509535548597

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("M3").Value = 6.153936333333334
$ws.Range("O3").Value = 0.5213583040808726
$ws.Range("P3").Value = 0.5213583040808725
$ws.Range("Q3").Value = 0.1448985335925556
$ws.Range("R3").Value = 1.304086802333
$ws.Range("S3").Value = 0.001467913184472214
$ws.Range("T3").Value = 0.001467913184472214

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.2006677815562353
$ws.Range("P4").Value = 0.2006677815562353
$ws.Range("Q4").Value = 0.05577060355455556
$ws.Range("R4").Value = 0.501935431991
$ws.Range("S4").Value = 0.0005649912544588448
$ws.Range("T4").Value = 0.0005649912544588446

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 3.281109666666667
$ws.Range("N5").Value = 9.843329000000001
$ws.Range("O5").Value = 0.2779739143628921
$ws.Range("P5").Value = 0.2779739143628921
$ws.Range("Q5").Value = 27.07801921461612
$ws.Range("R5").Value = 243.702172931545
$ws.Range("S5").Value = 0.2743173476571958
$ws.Range("T5").Value = 0.2743173476571958

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("M6").Value = 6.153936333333334
$ws.Range("O6").Value = 0.5213583040808726
$ws.Range("P6").Value = 0.5213583040808725
$ws.Range("Q6").Value = 50.78660063466057
$ws.Range("R6").Value = 457.0794057119451
$ws.Range("S6").Value = 0.5145001734508465
$ws.Range("T6").Value = 0.5145001734508463

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.2006677815562353
$ws.Range("P7").Value = 0.2006677815562353
$ws.Range("S7").Value = 0.1980281269302744
$ws.Range("T7").Value = 0.1980281269302744

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 3.281109666666667
$ws.Range("N8").Value = 9.843329000000001
$ws.Range("O8").Value = 0.2779739143628921
$ws.Range("P8").Value = 0.2779739143628921
$ws.Range("Q8").Value = 0.2836858354832222
$ws.Range("R8").Value = 2.553172519349
$ws.Range("S8").Value = 0.002873915752141414
$ws.Range("T8").Value = 0.002873915752141414

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("M9").Value = 6.153936333333334
$ws.Range("O9").Value = 0.5213583040808726
$ws.Range("P9").Value = 0.5213583040808725
$ws.Range("R9").Value = 4.788642480229
$ws.Range("S9").Value = 0.00539021744555385
$ws.Range("T9").Value = 0.005390217445553848

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.2006677815562353
$ws.Range("P10").Value = 0.2006677815562353
$ws.Range("S10").Value = 0.00207466337150204
$ws.Range("T10").Value = 0.002074663371502039
